# Update the "Training Dashboard" sheet for the new progress snapshot
# taken on 04-Nov-2025:
#   - column H (PERIOD TO EXPIRE) drops by 1 day for every data row
#   - column I (LAST UPDATE) moves from 03-Nov-2025 to 04-Nov-2025
#
# Data rows run from row 3 through row 23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$xlPasteValues = -4163

for ($row = 3; $row -le 23; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE
    # Writing the date-like text directly through Value/Value2 gets
    # auto-converted into a real date serial + date number format by
    # the engine. Route it through a formula -> copy -> paste-values
    # round trip instead, which keeps the cell a plain text value (as
    # it was before) with its original "General" number format/style.
    $iCell.Formula = '="04-Nov-2025"'
    $iCell.Copy()
    $iCell.PasteSpecial($xlPasteValues)
}
